$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 21:35"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1634790
$ws.Range("C4").Value = 13888
$ws.Range("D4").Value = 386081
$ws.Range("E4").Value = 1151525
$ws.Range("G4").Value = 830
$ws.Range("H4").Value = 97184

# --- Alemania (row 11) ---
$ws.Range("B11").Value = 179712
$ws.Range("C11").Value = 691
$ws.Range("E11").Value = 12361
$ws.Range("G11").Value = 42
$ws.Range("H11").Value = 8351

# --- India (row 14) ---
$ws.Range("B14").Value = 124792
$ws.Range("C14").Value = 6566
$ws.Range("D14").Value = 51820
$ws.Range("E14").Value = 69246

# --- Islas Caimanes moves up the ranking into row 164 (new data), pushing
#     Guyana / Bermudas / Camboya each down one row (rows 165-167). The
#     country previously in row 168 (Trinidad y Tobago) is unaffected. ---

# Row 164: Islas Caimanes (new top position, updated stats)
$ws.Range("A164").Value = "Islas Caimanes"
$ws.Range("B164").Value = 129
$ws.Range("C164").Value = 8
$ws.Range("D164").Value = 61
$ws.Range("E164").Value = 67
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1

# Row 165: Guyana (shifted down from old row 164)
$ws.Range("A165").Value = "Guyana"
$ws.Range("B165").Value = 127
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 57
$ws.Range("E165").Value = 60
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 10

# Row 166: Bermudas (shifted down from old row 165)
$ws.Range("A166").Value = "Bermudas"
$ws.Range("B166").Value = 125
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 80
$ws.Range("E166").Value = 36
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 9

# Row 167: Camboya (shifted down from old row 166)
$ws.Range("A167").Value = "Camboya"
$ws.Range("B167").Value = 123
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 122
$ws.Range("E167").Value = 1
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0
